$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the Resolving-Mac rows (11-13); this also removes the now-unused shared string
$ws.Range("A11:T13").Delete() | Out-Null

# Update numeric data for rows 2-10 (columns E:T)
$row2 = New-Object 'object[,]' 1,16
$row2[0,0] = 3.0
$row2[0,1] = 1.0
$row2[0,2] = 0.6599316666666667
$row2[0,3] = 1.979795
$row2[0,4] = 0.004637879740437423
$row2[0,5] = 0.004642763851664535
$row2[0,6] = 3.0
$row2[0,7] = 1.0
$row2[0,8] = 12.36292333333333
$row2[0,9] = 37.08877
$row2[0,10] = 0.918273862214392
$row2[0,11] = 0.930379446606803
$row2[0,12] = 8.158684600238889
$row2[0,13] = 73.42816140215
$row2[0,14] = 0.004258843741737355
$row2[0,15] = 0.00431953206303772
$ws.Range("E2:T2").Value2 = $row2

$row3 = New-Object 'object[,]' 1,16
$row3[0,0] = 3.0
$row3[0,1] = 1.0
$row3[0,2] = 0.6599316666666667
$row3[0,3] = 1.979795
$row3[0,4] = 0.004637879740437423
$row3[0,5] = 0.004642763851664535
$row3[0,6] = 3.0
$row3[0,7] = 1.0
$row3[0,8] = 0.574769
$row3[0,9] = 1.724307
$row3[0,10] = 0.04269179184247177
$row3[0,11] = 0.04325459680761149
$row3[0,12] = 0.3793082641183333
$row3[0,13] = 3.413774377065
$row3[0,14] = 0.0001979993964691715
$row3[0,15] = 0.0002008208784767028
$ws.Range("E3:T3").Value2 = $row3

$row4 = New-Object 'object[,]' 1,16
$row4[0,0] = 3.0
$row4[0,1] = 1.0
$row4[0,2] = 0.6599316666666667
$row4[0,3] = 1.979795
$row4[0,4] = 0.004637879740437423
$row4[0,5] = 0.004642763851664535
$row4[0,6] = 2.0
$row4[0,7] = 1.0
$row4[0,8] = 0.525528
$row4[0,9] = 1.051056
$row4[0,10] = 0.03903434594313629
$row4[0,11] = 0.02636595658558534
$row4[0,12] = 0.34681256892
$row4[0,13] = 2.08087541352
$row4[0,14] = 0.0001810366022308975
$row4[0,15] = 0.0001224109101501121
$ws.Range("E4:T4").Value2 = $row4

$row5 = New-Object 'object[,]' 1,16
$row5[0,0] = 3.0
$row5[0,1] = 1.0
$row5[0,2] = 141.1826756666667
$row5[0,3] = 423.548027
$row5[0,4] = 0.9922061695910651
$row5[0,5] = 0.9932510533663508
$row5[0,6] = 3.0
$row5[0,7] = 1.0
$row5[0,8] = 12.36292333333333
$row5[0,9] = 37.08877
$row5[0,10] = 0.918273862214392
$row5[0,11] = 0.930379446606803
$row5[0,12] = 1745.430595261865
$row5[0,13] = 15708.87535735679
$row5[0,14] = 0.9111169914633354
$row5[0,15] = 0.9241003653726096
$ws.Range("E5:T5").Value2 = $row5

$row6 = New-Object 'object[,]' 1,16
$row6[0,0] = 3.0
$row6[0,1] = 1.0
$row6[0,2] = 141.1826756666667
$row6[0,3] = 423.548027
$row6[0,4] = 0.9922061695910651
$row6[0,5] = 0.9932510533663508
$row6[0,6] = 3.0
$row6[0,7] = 1.0
$row6[0,8] = 0.574769
$row6[0,9] = 1.724307
$row6[0,10] = 0.04269179184247177
$row6[0,11] = 0.04325459680761149
$row6[0,12] = 81.14742531025432
$row6[0,13] = 730.326827792289
$row6[0,14] = 0.04235905925699799
$row6[0,15] = 0.04296267384209691
$ws.Range("E6:T6").Value2 = $row6

$row7 = New-Object 'object[,]' 1,16
$row7[0,0] = 3.0
$row7[0,1] = 1.0
$row7[0,2] = 141.1826756666667
$row7[0,3] = 423.548027
$row7[0,4] = 0.9922061695910651
$row7[0,5] = 0.9932510533663508
$row7[0,6] = 2.0
$row7[0,7] = 1.0
$row7[0,8] = 0.525528
$row7[0,9] = 1.051056
$row7[0,10] = 0.03903434594313629
$row7[0,11] = 0.02636595658558534
$row7[0,12] = 74.195449177752
$row7[0,13] = 445.172695066512
$row7[0,14] = 0.03873011887073179
$row7[0,15] = 0.02618801415164411
$ws.Range("E7:T7").Value2 = $row7

$row8 = New-Object 'object[,]' 1,16
$row8[0,0] = 2.0
$row8[0,1] = 1.0
$row8[0,2] = 0.4490655
$row8[0,3] = 0.898131
$row8[0,4] = 0.0031559506684976
$row8[0,5] = 0.00210618278198466
$row8[0,6] = 3.0
$row8[0,7] = 1.0
$row8[0,8] = 12.36292333333333
$row8[0,9] = 37.08877
$row8[0,10] = 0.918273862214392
$row8[0,11] = 0.930379446606803
$row8[0,12] = 5.551762348145
$row8[0,13] = 33.31057408887
$row8[0,14] = 0.002898027009319383
$row8[0,15] = 0.001959549171155665
$ws.Range("E8:T8").Value2 = $row8

$row9 = New-Object 'object[,]' 1,16
$row9[0,0] = 2.0
$row9[0,1] = 1.0
$row9[0,2] = 0.4490655
$row9[0,3] = 0.898131
$row9[0,4] = 0.0031559506684976
$row9[0,5] = 0.00210618278198466
$row9[0,6] = 3.0
$row9[0,7] = 1.0
$row9[0,8] = 0.574769
$row9[0,9] = 1.724307
$row9[0,10] = 0.04269179184247177
$row9[0,11] = 0.04325459680761149
$row9[0,12] = 0.2581089283695
$row9[0,13] = 1.548653570217
$row9[0,14] = 0.0001347331890046091
$row9[0,15] = 0.00009110208703787998
$ws.Range("E9:T9").Value2 = $row9

$row10 = New-Object 'object[,]' 1,16
$row10[0,0] = 2.0
$row10[0,1] = 1.0
$row10[0,2] = 0.4490655
$row10[0,3] = 0.898131
$row10[0,4] = 0.0031559506684976
$row10[0,5] = 0.00210618278198466
$row10[0,6] = 2.0
$row10[0,7] = 1.0
$row10[0,8] = 0.525528
$row10[0,9] = 1.051056
$row10[0,10] = 0.03903434594313629
$row10[0,11] = 0.02636595658558534
$row10[0,12] = 0.235996494084
$row10[0,13] = 0.943985976336
$row10[0,14] = 0.0001231904701736075
$row10[0,15] = 0.00005553152379111491
$ws.Range("E10:T10").Value2 = $row10

